$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1965575218200684
$ws.Range("E2").Value = 80.11389566444632
$ws.Range("F2").Value = 0.0034064004093979
$ws.Range("G2").Value = 0.002939700397567074
$ws.Range("H2").Value = 0.002566815087693311
$ws.Range("I2").Value = 0.002326657877359087
$ws.Range("J2").Value = 0.002286228920095569
$ws.Range("K2").Value = 0.002029458474717457
$ws.Range("L2").Value = 0.002029458474717457
$ws.Range("M2").Value = 0.002029458474717457
$ws.Range("N2").Value = 0.001876908969343045
$ws.Range("O2").Value = 0.00184622989131333
$ws.Range("P2").Value = 0.00184622989131333
$ws.Range("Q2").Value = 0.001773619890375423
$ws.Range("R2").Value = 0.001773619890375423
$ws.Range("S2").Value = 0.00169118143618738
$ws.Range("T2").Value = 0.001666108947535326
$ws.Range("U2").Value = 0.001640733283284011
$ws.Range("V2").Value = 0.001631100837795693
$ws.Range("W2").Value = 0.001588951921210529
$ws.Range("X2").Value = 0.001566405003710744
$ws.Range("Y2").Value = 0.001561674379423904
$ws.Range("C3").Value = 0.2621278762817383
$ws.Range("E3").Value = 80.98248494083418
$ws.Range("F3").Value = 0.003543165945709022
$ws.Range("G3").Value = 0.002938737812021651
$ws.Range("H3").Value = 0.002585375430321149
$ws.Range("I3").Value = 0.00246847138451027
$ws.Range("J3").Value = 0.002170182152740327
$ws.Range("K3").Value = 0.002058573563939306
$ws.Range("L3").Value = 0.001942078234194632
$ws.Range("M3").Value = 0.001942078234194632
$ws.Range("N3").Value = 0.001877640692208635
$ws.Range("O3").Value = 0.001756609446336681
$ws.Range("P3").Value = 0.001715321210203148
$ws.Range("Q3").Value = 0.001677190155027857
$ws.Range("R3").Value = 0.001677190155027857
$ws.Range("S3").Value = 0.00163185319085867
$ws.Range("T3").Value = 0.00160371113931548
$ws.Range("U3").Value = 0.00160371113931548
$ws.Range("V3").Value = 0.001595536182378044
$ws.Range("W3").Value = 0.001595536182378044
$ws.Range("X3").Value = 0.00158716833143148
$ws.Range("Y3").Value = 0.001578605944265773
$ws.Range("C4").Value = 0.2090551853179932
$ws.Range("E4").Value = 83.72578526828329
$ws.Range("F4").Value = 0.003690111078309416
$ws.Range("G4").Value = 0.002855100436841773
$ws.Range("H4").Value = 0.00224717225884831
$ws.Range("I4").Value = 0.002085178928853983
$ws.Range("J4").Value = 0.002085178928853983
$ws.Range("K4").Value = 0.002071198806504619
$ws.Range("L4").Value = 0.001988741344736346
$ws.Range("M4").Value = 0.00198174053723767
$ws.Range("N4").Value = 0.001797604959738166
$ws.Range("O4").Value = 0.001797604959738166
$ws.Range("P4").Value = 0.001760194688294998
$ws.Range("Q4").Value = 0.001760194688294998
$ws.Range("R4").Value = 0.001751304172548221
$ws.Range("S4").Value = 0.00172841520379968
$ws.Range("T4").Value = 0.001692372589404875
$ws.Range("U4").Value = 0.001674948440065305
$ws.Range("V4").Value = 0.001666155314114643
$ws.Range("W4").Value = 0.001660767118110849
$ws.Range("X4").Value = 0.001643487343924676
$ws.Range("Y4").Value = 0.001632081584177062
$ws.Range("C5").Value = 0.2100899219512939
$ws.Range("E5").Value = 80.19955851590021
$ws.Range("F5").Value = 0.003615637000805763
$ws.Range("G5").Value = 0.002853093296369826
$ws.Range("H5").Value = 0.002521850818783149
$ws.Range("I5").Value = 0.002269587982458932
$ws.Range("J5").Value = 0.002216245477540784
$ws.Range("K5").Value = 0.002099596092365649
$ws.Range("L5").Value = 0.001922077286967498
$ws.Range("M5").Value = 0.001854643122036786
$ws.Range("N5").Value = 0.001854643122036786
$ws.Range("O5").Value = 0.001835043641151794
$ws.Range("P5").Value = 0.001795015397663334
$ws.Range("Q5").Value = 0.001780761375098593
$ws.Range("R5").Value = 0.00171644628158132
$ws.Range("S5").Value = 0.00171644628158132
$ws.Range("T5").Value = 0.001651285661575623
$ws.Range("U5").Value = 0.001639365039344566
$ws.Range("V5").Value = 0.001601874947398181
$ws.Range("W5").Value = 0.001601874947398181
$ws.Range("X5").Value = 0.00156334422058285
$ws.Range("Y5").Value = 0.00156334422058285
$ws.Range("C6").Value = 0.2929625511169434
$ws.Range("E6").Value = 78.12934800561379
$ws.Range("F6").Value = 0.003465229645800299
$ws.Range("G6").Value = 0.002688582320078794
$ws.Range("H6").Value = 0.002526911795793723
$ws.Range("I6").Value = 0.002301782921082125
$ws.Range("J6").Value = 0.002055379795845787
$ws.Range("K6").Value = 0.002035526468536691
$ws.Range("L6").Value = 0.001913848647427106
$ws.Range("M6").Value = 0.001866845221754904
$ws.Range("N6").Value = 0.001833916652043055
$ws.Range("O6").Value = 0.001787789781182118
$ws.Range("P6").Value = 0.001700272239661316
$ws.Range("Q6").Value = 0.001669708583043434
$ws.Range("R6").Value = 0.001639516460510958
$ws.Range("S6").Value = 0.001619000711253376
$ws.Range("T6").Value = 0.00155587496782635
$ws.Range("U6").Value = 0.00155587496782635
$ws.Range("V6").Value = 0.00155587496782635
$ws.Range("W6").Value = 0.001544130462705822
$ws.Range("X6").Value = 0.00152873006226239
$ws.Range("Y6").Value = 0.001522989239875512
$ws.Range("C7").Value = 0.2398521900177002
$ws.Range("E7").Value = 79.17031733412296
$ws.Range("F7").Value = 0.003472056448682002
$ws.Range("G7").Value = 0.002883603473870912
$ws.Range("H7").Value = 0.002572335527272681
$ws.Range("I7").Value = 0.002255143558624631
$ws.Range("J7").Value = 0.002144623703995756
$ws.Range("K7").Value = 0.00202354122883451
$ws.Range("L7").Value = 0.001988364222944019
$ws.Range("M7").Value = 0.001906067559509758
$ws.Range("N7").Value = 0.001804498751647323
$ws.Range("O7").Value = 0.001803689794381137
$ws.Range("P7").Value = 0.001639768340158441
$ws.Range("Q7").Value = 0.001639768340158441
$ws.Range("R7").Value = 0.001639768340158441
$ws.Range("S7").Value = 0.001594879727396445
$ws.Range("T7").Value = 0.001594879727396445
$ws.Range("U7").Value = 0.001594879727396445
$ws.Range("V7").Value = 0.001579539492382363
$ws.Range("W7").Value = 0.001579539492382363
$ws.Range("X7").Value = 0.001563253973213223
$ws.Range("Y7").Value = 0.001543281039651519
$ws.Range("C8").Value = 0.2140746116638184
$ws.Range("E8").Value = 78.07068384169361
$ws.Range("G8").Value = 0.002932206378861229
$ws.Range("H8").Value = 0.002627136054231766
$ws.Range("I8").Value = 0.002325355616545497
$ws.Range("J8").Value = 0.00215020524883308
$ws.Range("K8").Value = 0.00215020524883308
$ws.Range("L8").Value = 0.00215020524883308
$ws.Range("M8").Value = 0.001959133325351119
$ws.Range("N8").Value = 0.001959133325351119
$ws.Range("O8").Value = 0.001695004734906318
$ws.Range("P8").Value = 0.001695004734906318
$ws.Range("Q8").Value = 0.001695004734906318
$ws.Range("R8").Value = 0.001695004734906318
$ws.Range("S8").Value = 0.001640810079011272
$ws.Range("T8").Value = 0.001614952861858144
$ws.Range("U8").Value = 0.00158010873233204
$ws.Range("V8").Value = 0.001576051803209413
$ws.Range("W8").Value = 0.001553221364669915
$ws.Range("X8").Value = 0.001529121290141527
$ws.Range("Y8").Value = 0.001521845688921902
$ws.Range("C9").Value = 0.21875
$ws.Range("E9").Value = 78.18120103215006
$ws.Range("F9").Value = 0.003744740124062058
$ws.Range("G9").Value = 0.002840132546236567
$ws.Range("H9").Value = 0.002204310088156319
$ws.Range("I9").Value = 0.002204310088156319
$ws.Range("J9").Value = 0.002153565057913645
$ws.Range("K9").Value = 0.002100089748053694
$ws.Range("L9").Value = 0.001952877674862464
$ws.Range("M9").Value = 0.001881658004841957
$ws.Range("N9").Value = 0.001822866030628819
$ws.Range("O9").Value = 0.001739706247139089
$ws.Range("P9").Value = 0.001711576716126026
$ws.Range("Q9").Value = 0.001711576716126026
$ws.Range("R9").Value = 0.001614114164554356
$ws.Range("S9").Value = 0.001571998599521291
$ws.Range("T9").Value = 0.001571998599521291
$ws.Range("U9").Value = 0.001571998599521291
$ws.Range("V9").Value = 0.001566797505886064
$ws.Range("W9").Value = 0.001544111956601399
$ws.Range("X9").Value = 0.001543437248945206
$ws.Range("Y9").Value = 0.001524000020119884
$ws.Range("C10").Value = 0.2479193210601807
$ws.Range("E10").Value = 77.48451050783478
$ws.Range("F10").Value = 0.003632288010573821
$ws.Range("G10").Value = 0.002961008011339377
$ws.Range("H10").Value = 0.002350267019669826
$ws.Range("I10").Value = 0.002236595090739905
$ws.Range("J10").Value = 0.002236595090739905
$ws.Range("K10").Value = 0.002099354709725124
$ws.Range("L10").Value = 0.002074970565844901
$ws.Range("M10").Value = 0.001956009252144088
$ws.Range("N10").Value = 0.00193864126303156
$ws.Range("O10").Value = 0.001866146243428765
$ws.Range("P10").Value = 0.001866014576791702
$ws.Range("Q10").Value = 0.001753418581980337
$ws.Range("R10").Value = 0.001753418581980337
$ws.Range("S10").Value = 0.001664772966093212
$ws.Range("T10").Value = 0.001664772966093212
$ws.Range("U10").Value = 0.001644558724383529
$ws.Range("V10").Value = 0.001596862942359768
$ws.Range("W10").Value = 0.001544399564283477
$ws.Range("X10").Value = 0.001512533860548866
$ws.Range("Y10").Value = 0.001510419308144927
$ws.Range("C11").Value = 0.2136814594268799
$ws.Range("E11").Value = 80.77882561079605
$ws.Range("F11").Value = 0.003731177597385016
$ws.Range("G11").Value = 0.002816787661299438
$ws.Range("H11").Value = 0.002579384315105313
$ws.Range("I11").Value = 0.002227024119271595
$ws.Range("J11").Value = 0.002227024119271595
$ws.Range("K11").Value = 0.002112213722177354
$ws.Range("L11").Value = 0.001910928604241963
$ws.Range("M11").Value = 0.001910928604241963
$ws.Range("N11").Value = 0.001910928604241963
$ws.Range("O11").Value = 0.001744282461701575
$ws.Range("P11").Value = 0.001688253714728241
$ws.Range("Q11").Value = 0.00166318533139617
$ws.Range("R11").Value = 0.00166318533139617
$ws.Range("S11").Value = 0.001620413131806298
$ws.Range("T11").Value = 0.001620413131806298
$ws.Range("U11").Value = 0.001620413131806298
$ws.Range("V11").Value = 0.00160969673835961
$ws.Range("W11").Value = 0.001599020981138466
$ws.Range("X11").Value = 0.001599020981138466
$ws.Range("Y11").Value = 0.001574635976818636
